# Weekly update: insert two new Choclo price records into the dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record at row 548 (shifts existing rows 548+ down by one) ---
$ws.Rows.Item(548).Insert()

$ws.Range("A548").Value = 8
$ws.Range("B548").Value = "Terminal La Palmera de La Serena"
$ws.Range("C548").Value = "Coquimbo"
$ws.Range("D548").Value = 44748
$ws.Range("E548").Value = 4
$ws.Range("F548").Value = 100112024
$ws.Range("G548").Value = "Choclo"
$ws.Range("H548").Value = "Dulce o Americano"
$ws.Range("I548").Value = "Primera"
$ws.Range("J548").Value = 500
$ws.Range("K548").Value = 43000
$ws.Range("L548").Value = 44000
$ws.Range("M548").Value = 43500
$ws.Range("N548").Value = "`$/malla 70 unidades"
$ws.Range("O548").Value = "Región de Arica y Parinacota"
$ws.Range("P548").Value = 621
$ws.Range("Q548").Value = 70
$ws.Range("R548").Value = "Hortaliza"

# --- Insert new record at row 605 (after the first insert has already shifted
#     rows down, this inserts a second new record and pushes the remainder
#     further down) ---
$ws.Rows.Item(605).Insert()

$ws.Range("A605").Value = 8
$ws.Range("B605").Value = "Terminal La Palmera de La Serena"
$ws.Range("C605").Value = "Coquimbo"
$ws.Range("D605").Value = 44747
$ws.Range("E605").Value = 4
$ws.Range("F605").Value = 100112024
$ws.Range("G605").Value = "Choclo"
$ws.Range("H605").Value = "Dulce o Americano"
$ws.Range("I605").Value = "Primera"
$ws.Range("J605").Value = 400
$ws.Range("K605").Value = 43000
$ws.Range("L605").Value = 44000
$ws.Range("M605").Value = 43500
$ws.Range("N605").Value = "`$/malla 70 unidades"
$ws.Range("O605").Value = "Región de Arica y Parinacota"
$ws.Range("P605").Value = 621
$ws.Range("Q605").Value = 70
$ws.Range("R605").Value = "Hortaliza"

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
